$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell E1 with the same style/formatting as the other headers
$ws.Range("E1").Value = "p-value"
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)

# Update row 2 values
$ws.Range("B2").Value = 0.4536790578878014
$ws.Range("C2").Value = 1.0086
$ws.Range("D2").Value = 0.006360000000000032
$ws.Range("E2").Value = 0.03487180888169839

# Update row 3 values
$ws.Range("B3").Value = 15.78029461131034
$ws.Range("C3").Value = 1.217369730421838
$ws.Range("D3").Value = 0.129969
$ws.Range("E3").Value = 0

# Update row 4 values
$ws.Range("B4").Value = -19.1220916628652
$ws.Range("C4").Value = 64.86818
$ws.Range("D4").Value = 0.1834885000000001
$ws.Range("E4").Value = 0

# Update row 5 values
$ws.Range("B5").Value = -5.171459538912004
$ws.Range("C5").Value = 91.60469999999999
$ws.Range("D5").Value = 0.06283219999999998
$ws.Range("E5").Value = 0
